$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 31 de Julio de 2020 a las 04:53"

# --- 2. Peru (row 10): refreshed case counts ---
$ws.Range("B10").Value = 407492
$ws.Range("D10").Value = 283915
$ws.Range("E10").Value = 104556
$ws.Range("H10").Value = 19021

# --- 3. Bolivia (row 35): refreshed case counts ---
$ws.Range("B35").Value = 75234
$ws.Range("C35").Value = 1700
$ws.Range("D35").Value = 23305
$ws.Range("E35").Value = 49035
$ws.Range("G35").Value = 86
$ws.Range("H35").Value = 2894

# --- 4. Belgica is inserted into the sorted table just after Ucrania (row 37). ---
#     Republica Dominicana and Bielorrusia (previously rows 38-39) each shift
#     down one row, and Belgica's fresh numbers land on row 38.
$ws.Range("A38").Value = "Belgica"
$ws.Range("B38").Value = 68006
$ws.Range("C38").Value = 671
$ws.Range("D38").Value = 17513
$ws.Range("E38").Value = 40653
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 4
$ws.Range("H38").Value = 9840

$ws.Range("A39").Value = "Republica Dominicana"
$ws.Range("B39").Value = 67915
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 35302
$ws.Range("E39").Value = 31467
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 1146

$ws.Range("A40").Value = "Bielorrusia"
$ws.Range("B40").Value = 67665
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 61765
$ws.Range("E40").Value = 5347
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 553

# --- 5. Honduras (row 51): refreshed case counts ---
$ws.Range("B51").Value = 41426
$ws.Range("C51").Value = 482
$ws.Range("D51").Value = 5443
$ws.Range("E51").Value = 34671
$ws.Range("G51").Value = 53
$ws.Range("H51").Value = 1312

# --- 6. Haiti (row 91): refreshed case counts ---
$ws.Range("B91").Value = 7412
$ws.Range("C91").Value = 34
$ws.Range("E91").Value = 2784
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = 161

# --- 7. Nueva Zelanda (row 137): refreshed case counts ---
$ws.Range("D137").Value = 1518
$ws.Range("E137").Value = 20

# --- 8. Camboya (row 176): refreshed case counts ---
$ws.Range("D176").Value = 164
$ws.Range("E176").Value = 70
